$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2
Set-TextValue "D2" "28.297.86"
$ws.Range("E2").Value = "  +2.82%  "

# Row 3
Set-TextValue "D3" "1.871.45"
$ws.Range("E3").Value = "  +1.57%  "

# Row 4
Set-TextValue "D4" "1.001"
$ws.Range("E4").Value = "  -0.44%  "

# Row 5
Set-TextValue "D5" "336.71"
$ws.Range("E5").Value = "  +0.96%  "

# Row 6
Set-TextValue "D6" "0.9998"
$ws.Range("E6").Value = "  -0.53%  "

# Row 7
Set-TextValue "D7" "0.4697"
$ws.Range("E7").Value = "  +1.43%  "

# Row 8
Set-TextValue "D8" "0.3925"
$ws.Range("E8").Value = "  +2.02%  "

# Row 9
Set-TextValue "D9" "47.12"
$ws.Range("E9").Value = "  +2.39%  "

# Row 10
Set-TextValue "D10" "0.07981"
$ws.Range("E10").Value = "  +0.94%  "

# Row 11
Set-TextValue "D11" "1.011"
$ws.Range("E11").Value = "  +1.88%  "

# Row 12
Set-TextValue "D12" "21.75"
$ws.Range("E12").Value = "  +1.23%  "

# Row 13
Set-TextValue "D13" "5.994"
$ws.Range("E13").Value = "  +1.16%  "

# Row 16
Set-TextValue "D16" "91.10"
$ws.Range("E16").Value = "  +2.65%  "

# Row 17
Set-TextValue "D17" "1.000"
$ws.Range("E17").Value = "  -0.58%  "

# Row 18
Set-TextValue "D18" "0.00001042"
$ws.Range("E18").Value = "  +0.70%  "

# Row 19
Set-TextValue "D19" "0.06579"
$ws.Range("E19").Value = "  -1.47%  "

# Row 20
Set-TextValue "D20" "17.67"
$ws.Range("E20").Value = "  +3.49%  "

# Row 21
Set-TextValue "D21" "0.9993"
$ws.Range("E21").Value = "  -0.61%  "

# Row 22
Set-TextValue "D22" "28.320.37"
$ws.Range("E22").Value = "  +2.87%  "

# Row 23
Set-TextValue "D23" "5.457"
$ws.Range("E23").Value = "  +1.46%  "

# Row 24
Set-TextValue "D24" "11.05"
$ws.Range("E24").Value = "  +1.45%  "

# Row 25
Set-TextValue "D25" "2.288"
$ws.Range("E25").Value = "  -1.01%  "

# Row 26
Set-TextValue "D26" "2.073.48"

# Row 27
Set-TextValue "D27" "158.79"
$ws.Range("E27").Value = "  -0.10%  "

# Row 28
Set-TextValue "D28" "19.95"
$ws.Range("E28").Value = "  +2.42%  "

# Row 29
Set-TextValue "D29" "2.154"
$ws.Range("E29").Value = "  +2.54%  "

# Row 30
Set-TextValue "D30" "5.526"
$ws.Range("E30").Value = "  +2.40%  "

# Row 31
Set-TextValue "D31" "119.88"
$ws.Range("E31").Value = "  +0.05%  "

# Row 32
Set-TextValue "D32" "0.9793"
$ws.Range("E32").Value = "  +0.52%  "

# Row 33
Set-TextValue "D33" "0.09473"
$ws.Range("E33").Value = "  +0.97%  "

# Row 34
Set-TextValue "D34" "3.572"
$ws.Range("E34").Value = "  -0.60%  "

# Row 35
Set-TextValue "D35" "5.366"
$ws.Range("E35").Value = "  +1.76%  "

# Row 36
Set-TextValue "D36" "1.376"
$ws.Range("E36").Value = "  +2.70%  "

# Row 37
Set-TextValue "D37" "0.02265"
$ws.Range("E37").Value = "  +1.71%  "

# Row 38
Set-TextValue "D38" "0.06103"
$ws.Range("E38").Value = "  +1.30%  "

# Row 39
Set-TextValue "D39" "8.409"
$ws.Range("E39").Value = "  +1.20%  "

# Row 40
Set-TextValue "D40" "1.179"
$ws.Range("E40").Value = "  -0.37%  "

# Row 41
Set-TextValue "D41" "0.5980"
$ws.Range("E41").Value = "  +1.71%  "

# Row 42
Set-TextValue "D42" "0.9985"
$ws.Range("E42").Value = "  -0.65%  "

# Row 43
Set-TextValue "D43" "0.1882"
$ws.Range("E43").Value = "  +0.91%  "

# Row 44
Set-TextValue "D44" "10.41"
$ws.Range("E44").Value = "  +1.10%  "

# Row 45
Set-TextValue "D45" "1.283"
$ws.Range("E45").Value = "  +3.51%  "

# Row 46
Set-TextValue "D46" "0.5595"
$ws.Range("E46").Value = "  +0.37%  "

# Row 47
Set-TextValue "D47" "12.27"
$ws.Range("E47").Value = "  +1.09%  "

# Row 48
Set-TextValue "D48" "1.970"
$ws.Range("E48").Value = "  +3.59%  "

# Row 49
Set-TextValue "D49" "0.06896"
$ws.Range("E49").Value = "  +3.00%  "

# Row 50
Set-TextValue "D50" "110.99"
$ws.Range("E50").Value = "  +0.01%  "

# Row 51
Set-TextValue "D51" "1.980"
$ws.Range("E51").Value = "  +11.75%  "

# Row 14 <-> Row 15 swap (Chainlink / WrappedEther)
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D14" "1.857.07"
$ws.Range("E14").Value = "  +0.04%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D15" "7.268"
$ws.Range("E15").Value = "  +2.20%  "
